# Fix the ticker typo "APPL" -> "AAPL" across all worksheets.
# (The ticker for Apple Inc. is AAPL, not APPL.)

$wb = $excel.ActiveWorkbook

$deposits = $wb.Worksheets.Item("deposits")
$dividends = $wb.Worksheets.Item("dividends")
$sales = $wb.Worksheets.Item("sales")

# deposits!B6 : APPL -> AAPL
$deposits.Range("B6").Value = "AAPL"

# dividends!B3 : APPL -> AAPL
$dividends.Range("B3").Value = "AAPL"

# sales!B6 and sales!B7 : APPL -> AAPL
$sales.Range("B6").Value = "AAPL"
$sales.Range("B7").Value = "AAPL"
